$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsARM = $wb.Worksheets.Item("ARM")
$wsBSM = $wb.Worksheets.Item("BSM")
$wsCRP = $wb.Worksheets.Item("CRP")
$wsCUL = $wb.Worksheets.Item("CUL")
$wsGSM = $wb.Worksheets.Item("GSM")
$wsLTW = $wb.Worksheets.Item("LTW")
$wsWVR = $wb.Worksheets.Item("WVR")

# Row 17: One for the Road / Potion
$wsALC.Range("H17").Value = 643855.6
$wsALC.Range("J17").Value = 643855.6
$wsALC.Range("L17").Value = 1931566.8
$wsALC.Range("N17").Value = -1931902.8

# Row 19: Unbreak My Heart / Roof Tile
$wsALC.Range("H19").Value = 406.3125
$wsALC.Range("I19").Value = 392.2
$wsALC.Range("J19").Value = 412.72726
$wsALC.Range("K19").Value = 392.2
$wsALC.Range("L19").Value = 412.72726
$wsALC.Range("M19").Value = -217.2
$wsALC.Range("N19").Value = -762.72726

# Row 64: Forged from the Void / Void Glue
$wsALC.Range("H64").Value = 2773.9722
$wsALC.Range("I64").Value = 2690.2
$wsALC.Range("J64").Value = 2878.6875
$wsALC.Range("K64").Value = 2690.2
$wsALC.Range("L64").Value = 2878.6875
$wsALC.Range("M64").Value = -2442.2
$wsALC.Range("N64").Value = -3374.6875

# Row 67: Dodging the Draft (L) / Void Glue
$wsALC.Range("H67").Value = 2773.9722
$wsALC.Range("I67").Value = 2690.2
$wsALC.Range("J67").Value = 2878.6875
$wsALC.Range("K67").Value = 2690.2
$wsALC.Range("L67").Value = 2878.6875
$wsALC.Range("M67").Value = -1832.2
$wsALC.Range("N67").Value = -4594.6875

# Row 103: Let Loose the Juice / Persimmon Tannin
$wsALC.Range("H103").Value = 334.66666
$wsALC.Range("I103").Value = 334.66666
$wsALC.Range("J103").Value = 0
$wsALC.Range("K103").Value = 1003.99998
$wsALC.Range("L103").Value = 0
$wsALC.Range("M103").Value = -417.9999799999999
$wsALC.Range("N103").ClearContents()

# Row 32: Ingot We Trust / Steel Ingot
$wsARM.Range("H32").Value = 23177.383
$wsARM.Range("I32").Value = 5450.811
$wsARM.Range("J32").Value = 88765.7
$wsARM.Range("K32").Value = 5450.811
$wsARM.Range("L32").Value = 88765.7
$wsARM.Range("M32").Value = -5163.811
$wsARM.Range("N32").Value = -89339.7

# Row 102: Smells of Rich Tama-hagane / Tama-hagane Ingot
$wsARM.Range("H102").Value = 0
$wsARM.Range("I102").Value = 0
$wsARM.Range("K102").Value = 0
$wsARM.Range("M102").ClearContents()

# Row 105: Ingot to Wing It / Molybdenum Ingot
$wsBSM.Range("H105").Value = 2777.8235
$wsBSM.Range("I105").Value = 3282.9092
$wsBSM.Range("J105").Value = 1851.8334
$wsBSM.Range("K105").Value = 3282.9092
$wsBSM.Range("L105").Value = 1851.8334
$wsBSM.Range("M105").Value = -1535.9092
$wsBSM.Range("N105").Value = -5345.8334

# Row 99: O Pine / Pine Lumber
$wsCRP.Range("H99").Value = 1681.125
$wsCRP.Range("I99").Value = 1273.6666
$wsCRP.Range("J99").Value = 1925.6
$wsCRP.Range("K99").Value = 1273.6666
$wsCRP.Range("L99").Value = 1925.6
$wsCRP.Range("M99").Value = 224.3334
$wsCRP.Range("N99").Value = -4921.6

# Row 126: A Better Conductor / Red Pine Lumber
$wsCRP.Range("H126").Value = 1681.125
$wsCRP.Range("I126").Value = 1273.6666
$wsCRP.Range("J126").Value = 1925.6
$wsCRP.Range("K126").Value = 3820.9998
$wsCRP.Range("L126").Value = 5776.799999999999
$wsCRP.Range("M126").Value = -1350.9998
$wsCRP.Range("N126").Value = -10716.8

# Row 20: Omelette's Be Friends / Dodo Omelette
$wsCUL.Range("H20").Value = 800
$wsCUL.Range("I20").Value = 800
$wsCUL.Range("J20").Value = 0
$wsCUL.Range("K20").Value = 2400
$wsCUL.Range("L20").Value = 0
$wsCUL.Range("M20").Value = -2173
$wsCUL.Range("N20").ClearContents()

# Row 22: A Total Nut Job / Walnut Bread
$wsCUL.Range("H22").Value = 1050
$wsCUL.Range("I22").Value = 0
$wsCUL.Range("J22").Value = 1050
$wsCUL.Range("K22").Value = 0
$wsCUL.Range("L22").Value = 3150
$wsCUL.Range("M22").ClearContents()
$wsCUL.Range("N22").Value = -3488

# Row 24: Rustic Repast / Chicken and Mushrooms
$wsCUL.Range("H24").Value = 1045
$wsCUL.Range("I24").Value = 0
$wsCUL.Range("J24").Value = 1045
$wsCUL.Range("K24").Value = 0
$wsCUL.Range("L24").Value = 3135
$wsCUL.Range("M24").ClearContents()
$wsCUL.Range("N24").Value = -3595

# Row 25: Flakes for Friends / Apple Tart
$wsCUL.Range("H25").Value = 450
$wsCUL.Range("I25").Value = 450
$wsCUL.Range("K25").Value = 1350
$wsCUL.Range("M25").Value = -1181

# Row 27: Brain Food / Walnut Bread
$wsCUL.Range("H27").Value = 1050
$wsCUL.Range("I27").Value = 0
$wsCUL.Range("J27").Value = 1050
$wsCUL.Range("K27").Value = 0
$wsCUL.Range("L27").Value = 3150
$wsCUL.Range("M27").ClearContents()
$wsCUL.Range("N27").Value = -3354

# Row 30: Picnic Panic / Apple Tart
$wsCUL.Range("H30").Value = 450
$wsCUL.Range("I30").Value = 450
$wsCUL.Range("K30").Value = 1350
$wsCUL.Range("M30").Value = -1248

# Row 33: Cooking with Gas / Chicken Stock
$wsCUL.Range("H33").Value = 232.4
$wsCUL.Range("I33").Value = 6
$wsCUL.Range("J33").Value = 383.33334
$wsCUL.Range("K33").Value = 36
$wsCUL.Range("L33").Value = 2300.00004
$wsCUL.Range("M33").Value = 247
$wsCUL.Range("N33").Value = -2866.00004

# Row 34: Fever Pitch / Chamomile Tea
$wsCUL.Range("H34").Value = 2733.3333
$wsCUL.Range("I34").Value = 5000
$wsCUL.Range("J34").Value = 2280
$wsCUL.Range("K34").Value = 15000
$wsCUL.Range("L34").Value = 6840
$wsCUL.Range("M34").Value = -14916
$wsCUL.Range("N34").Value = -7008

# Row 49: Leek Soup for the Soul / Cawl Cennin
$wsCUL.Range("H49").Value = 10000
$wsCUL.Range("J49").Value = 10000
$wsCUL.Range("L49").Value = 30000
$wsCUL.Range("N49").Value = -30312

# Row 54: Good Eats in Ishgard / Salt Cod Puffs
$wsCUL.Range("H54").Value = 5988
$wsCUL.Range("I54").Value = 2599.5
$wsCUL.Range("J54").Value = 7117.5
$wsCUL.Range("K54").Value = 7798.5
$wsCUL.Range("L54").Value = 21352.5
$wsCUL.Range("M54").Value = -7239.5
$wsCUL.Range("N54").Value = -22470.5

# Row 55: Pagan Pastries / Pastry Fish
$wsCUL.Range("H55").Value = 2907.2666
$wsCUL.Range("I55").Value = 1004
$wsCUL.Range("J55").Value = 3043.2144
$wsCUL.Range("K55").Value = 3012
$wsCUL.Range("L55").Value = 9129.643199999999
$wsCUL.Range("N55").Value = -9483.643199999999
$wsCUL.Range("M55").Value = -2835

# Row 59: Comfort Me with Mushrooms / Buttons in a Blanket
$wsCUL.Range("H59").Value = 20000
$wsCUL.Range("I59").Value = 0
$wsCUL.Range("J59").Value = 20000
$wsCUL.Range("K59").Value = 0
$wsCUL.Range("L59").Value = 60000
$wsCUL.Range("M59").ClearContents()
$wsCUL.Range("N59").Value = -61080

# Row 60: Drinking to Your Health / Mulled Tea
$wsCUL.Range("H60").Value = 1345.125
$wsCUL.Range("I60").Value = 191.25
$wsCUL.Range("J60").Value = 2499
$wsCUL.Range("K60").Value = 573.75
$wsCUL.Range("L60").Value = 7497
$wsCUL.Range("M60").Value = -322.75
$wsCUL.Range("N60").Value = -7999

# Row 61: Red Letter Day / Rolanberry Lassi
$wsCUL.Range("H61").Value = 286
$wsCUL.Range("I61").Value = 300
$wsCUL.Range("J61").Value = 284
$wsCUL.Range("K61").Value = 900
$wsCUL.Range("L61").Value = 852
$wsCUL.Range("M61").Value = -685
$wsCUL.Range("N61").Value = -1282

# Row 68: Such a Butter Face / Fermented Butter
$wsCUL.Range("H68").Value = 11480.333
$wsCUL.Range("I68").Value = 388.57144
$wsCUL.Range("K68").Value = 1165.71432
$wsCUL.Range("M68").Value = -354.71432

# Row 69: Loving That Muffin Top / Ishgardian Muffin
$wsCUL.Range("H69").Value = 733.3333
$wsCUL.Range("J69").Value = 733.3333
$wsCUL.Range("L69").Value = 2199.9999
$wsCUL.Range("N69").Value = -3821.9999

# Row 71: No Margarine of Error (L) / Fermented Butter
$wsCUL.Range("H71").Value = 11480.333
$wsCUL.Range("I71").Value = 388.57144
$wsCUL.Range("K71").Value = 3497.14296
$wsCUL.Range("M71").Value = 558.8570399999999

# Row 72: Muffin of the Morn (L) / Ishgardian Muffin
$wsCUL.Range("H72").Value = 733.3333
$wsCUL.Range("J72").Value = 733.3333
$wsCUL.Range("L72").Value = 6599.9997
$wsCUL.Range("N72").Value = -14711.9997

# Row 76: Old Victories, New Tastes / Dhalmel Fricassee
$wsCUL.Range("H76").Value = 10400
$wsCUL.Range("I76").Value = 1000
$wsCUL.Range("J76").Value = 16666.666
$wsCUL.Range("K76").Value = 3000
$wsCUL.Range("L76").Value = 49999.99800000001
$wsCUL.Range("N76").Value = -50765.99800000001
$wsCUL.Range("M76").Value = -2617

# Row 79: The Eats of Authenticity (L) / Dhalmel Fricassee
$wsCUL.Range("H79").Value = 10400
$wsCUL.Range("I79").Value = 1000
$wsCUL.Range("J79").Value = 16666.666
$wsCUL.Range("K79").Value = 3000
$wsCUL.Range("L79").Value = 49999.99800000001
$wsCUL.Range("N79").Value = -52651.99800000001
$wsCUL.Range("M79").Value = -1674

# Row 80: Saucy for a Suitor / Hollandaise Sauce
$wsCUL.Range("H80").Value = 7523
$wsCUL.Range("J80").Value = 7773.6665
$wsCUL.Range("L80").Value = 23320.9995
$wsCUL.Range("N80").Value = -25192.9995

# Row 81: It Goes Down Smoothly / Frozen Spirits
$wsCUL.Range("H81").Value = 6500
$wsCUL.Range("I81").Value = 0
$wsCUL.Range("J81").Value = 6500
$wsCUL.Range("K81").Value = 0
$wsCUL.Range("L81").Value = 19500
$wsCUL.Range("M81").ClearContents()
$wsCUL.Range("N81").Value = -21746

# Row 82: Persuasion of a Higher Power / Baked Pipira Pira
$wsCUL.Range("H82").Value = 16763
$wsCUL.Range("J82").Value = 20203.75
$wsCUL.Range("L82").Value = 60611.25
$wsCUL.Range("N82").Value = -61423.25

# Row 83: Saved by the Sauce (L) / Hollandaise Sauce
$wsCUL.Range("H83").Value = 7523
$wsCUL.Range("J83").Value = 7773.6665
$wsCUL.Range("L83").Value = 69962.9985
$wsCUL.Range("N83").Value = -79322.9985

# Row 84: Quenching the Flame (L) / Frozen Spirits
$wsCUL.Range("H84").Value = 6500
$wsCUL.Range("I84").Value = 0
$wsCUL.Range("J84").Value = 6500
$wsCUL.Range("K84").Value = 0
$wsCUL.Range("L84").Value = 58500
$wsCUL.Range("M84").ClearContents()
$wsCUL.Range("N84").Value = -69732

# Row 85: Loaves and Fishes (L) / Baked Pipira Pira
$wsCUL.Range("H85").Value = 16763
$wsCUL.Range("J85").Value = 20203.75
$wsCUL.Range("L85").Value = 60611.25
$wsCUL.Range("N85").Value = -63419.25

# Row 93: Loquacious / Loquat Juice
$wsCUL.Range("H93").Value = 5037.467
$wsCUL.Range("J93").Value = 5037.467
$wsCUL.Range("L93").Value = 15112.401
$wsCUL.Range("N93").Value = -18856.401

# Row 97: The Frier Never Lies / Cottonseed Oil
$wsCUL.Range("H97").Value = 687.8461
$wsCUL.Range("I97").Value = 180.25
$wsCUL.Range("K97").Value = 540.75
$wsCUL.Range("M97").Value = -44.75

# Row 98: Sweet Kiss of Death / Rice Vinegar
$wsCUL.Range("H98").Value = 770129.4
$wsCUL.Range("I98").Value = 764.7778
$wsCUL.Range("J98").Value = 2501199.8
$wsCUL.Range("K98").Value = 2294.3334
$wsCUL.Range("L98").Value = 7503599.399999999
$wsCUL.Range("M98").Value = -796.3334
$wsCUL.Range("N98").Value = -7506595.399999999

# Row 101: No Othard Choice / Egg Foo Young
$wsCUL.Range("H101").Value = 11014.5
$wsCUL.Range("J101").Value = 11014.5
$wsCUL.Range("L101").Value = 33043.5
$wsCUL.Range("N101").Value = -37911.5

# Row 102: Persimmony Snicket / Persimmon Pudding
$wsCUL.Range("H102").Value = 3200.4285
$wsCUL.Range("I102").Value = 4006
$wsCUL.Range("J102").Value = 3138.4614
$wsCUL.Range("K102").Value = 12018
$wsCUL.Range("L102").Value = 9415.3842
$wsCUL.Range("M102").Value = -9584
$wsCUL.Range("N102").Value = -14283.3842

# Row 104: Fits to a Tea / Doman Tea
$wsCUL.Range("H104").Value = 4082.25
$wsCUL.Range("J104").Value = 4082.25
$wsCUL.Range("L104").Value = 12246.75
$wsCUL.Range("N104").Value = -17488.75

# Row 126: Gold Rush Order / Phrygian Gold Ingot
$wsGSM.Range("H126").Value = 1439.1
$wsGSM.Range("I126").Value = 1308.6666
$wsGSM.Range("J126").Value = 1495
$wsGSM.Range("K126").Value = 3925.9998
$wsGSM.Range("L126").Value = 4485
$wsGSM.Range("M126").Value = -1455.9998
$wsGSM.Range("N126").Value = -9425

# Row 82: Trainin' the Neck / Dragon Leather
$wsLTW.Range("H82").Value = 2064.1482
$wsLTW.Range("I82").Value = 953.55554
$wsLTW.Range("J82").Value = 2619.4443
$wsLTW.Range("K82").Value = 953.55554
$wsLTW.Range("L82").Value = 2619.4443
$wsLTW.Range("M82").Value = -592.55554
$wsLTW.Range("N82").Value = -3341.4443

# Row 85: Training Is Only Skintight (L) / Dragon Leather
$wsLTW.Range("H85").Value = 2064.1482
$wsLTW.Range("I85").Value = 953.55554
$wsLTW.Range("J85").Value = 2619.4443
$wsLTW.Range("K85").Value = 953.55554
$wsLTW.Range("L85").Value = 2619.4443
$wsLTW.Range("M85").Value = 294.44446
$wsLTW.Range("N85").Value = -5115.4443

# Row 126: A Polished Purchase / Snow Linen
$wsWVR.Range("H126").Value = 589110.4399999999
$wsWVR.Range("I126").Value = 714899.4399999999
$wsWVR.Range("J126").Value = 2095
$wsWVR.Range("K126").Value = 2144698.32
$wsWVR.Range("L126").Value = 6285
$wsWVR.Range("M126").Value = -2142228.32
$wsWVR.Range("N126").Value = -11225

# Row 132: Comfy Cabins / Snow Cotton Cloth
$wsWVR.Range("H132").Value = 1038.8043
$wsWVR.Range("I132").Value = 692.72
$wsWVR.Range("J132").Value = 1450.8096
$wsWVR.Range("K132").Value = 2078.16
$wsWVR.Range("L132").Value = 4352.4288
$wsWVR.Range("M132").Value = 451.8400000000001
$wsWVR.Range("N132").Value = -9412.4288
